$d = $word.ActiveDocument

# Delete before[46:52] (no replacement)
$r = $d.Range($d.Paragraphs.Item(46).Range.Start, $d.Paragraphs.Item(52).Range.End)
$r.Delete()

# Replace before[36:44] (9 para) with 5 new paragraph(s)
$r = $d.Range($d.Paragraphs.Item(36).Range.Start, $d.Paragraphs.Item(44).Range.End)
$r.Delete()
$cur = $d.Paragraphs.Item(35)
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Impact"
$cur.Style = "Heading3"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Discovered systematic race coding errors affecting all Black and Asian-American voters"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Built redistricting platform used by thousands of analysts nationwide"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%"
$cur.Style = "Normal"

# Replace before[31:34] (4 para) with 12 new paragraph(s)
$r = $d.Range($d.Paragraphs.Item(31).Range.Start, $d.Paragraphs.Item(34).Range.End)
$r.Delete()
$cur = $d.Paragraphs.Item(30)
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "National Redistricting Platform (2020 - 2021)"
$cur.Style = "Heading3"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Cloud-based GeoDjango platform for redistricting analysis with real-time collaborative editing and Census integration, used by thousands of analysts nationwide"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Technologies: GeoDjango, PostGIS, AWS, Docker, React, Python"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Impact: Reduced mapping costs by 73.5%, saving organizations `$4.7M in operational expenses"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "FLEEM Political Polling System (2010 - 2012)"
$cur.Style = "Heading3"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Completely self-built IVR system using Twilio API that contacted tens of thousands of voters daily, replicated call center functionality to performance parity"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Technologies: Twilio API, Python, Django, PostgreSQL, JavaScript"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Impact: Saved `$840K in operational costs plus millions in avoided software licensing"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Geospatial Demographic Classification System (2013 - 2016)"
$cur.Style = "Heading3"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Machine learning platform that discovered systematic coding errors and improved demographic classification accuracy from 23% to 64%"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Technologies: Python, Scikit-learn, PostGIS, GeoPandas, TensorFlow"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Impact: Corrected demographic data affecting all Black and Asian-American voters nationwide"
$cur.Style = "Normal"

# Replace before[15:29] (15 para) with 3 new paragraph(s)
$r = $d.Range($d.Paragraphs.Item(15).Range.Start, $d.Paragraphs.Item(29).Range.End)
$r.Delete()
$cur = $d.Paragraphs.Item(14)
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Harmonized data from 20+ polling firms with incompatible methodologies and encoding systems"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+"
$cur.Style = "Normal"

# Replace before[8:13] (6 para) with 36 new paragraph(s)
$r = $d.Range($d.Paragraphs.Item(8).Range.Start, $d.Paragraphs.Item(13).Range.End)
$r.Delete()
$cur = $d.Paragraphs.Item(7)
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Partner - Siege Analytics (Austin, TX) | 2005 - Present"
$cur.Style = "Heading3"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Data, Technology and Strategy Consulting"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Data Products Manager - Helm/Murmuration (Austin, TX) | June 2021 - May 2023"
$cur.Style = "Heading3"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Civic Graph & Civic Pulse Director"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Conceived, architected and built Civic Graph multi-tenant data warehouse processing government data from Census, Bureau of Labor Statistics, National Council of Educational Statistics"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Built multi-dimensional data warehouse measuring socio-economic changes in America at every level across attitudinal, behavioral, demographic, economic and geographical dimensions"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Managed engineering teams of 7-11 professionals while setting technical direction for data architecture"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Analytics Supervisor - GSD&M (Austin, TX) | November 2019 - June 2020"
$cur.Style = "Heading3"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Big Data Engineering Transformation"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Transformed small data team into big data engineering team, scaling from laptop datasets to Hadoop Clusters and Hive on AWS"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Managed accounts including United States Air Force, Southwest Airlines/Chase and Indeed"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Rewrote mission and offerings of department and drafted integration plan with strategy team"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Software Engineer - Mautinoa Technologies (Austin, TX) | August 2016 - February 2018"
$cur.Style = "Heading3"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "SimCrisis Product Owner/Engineer"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Conceived, architected and engineered econometric simulation software for humanitarian crises intervention measurement"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Built SimCrisis GeoDjango web application using multi-agent modeling to create econometric simulations of crisis economies"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Designed modular application accepting rules extensions for ethnic strife, different crises/disasters, supply failures"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Senior Analyst - Myers Research (Austin, TX) | August 2012 - February 2014"
$cur.Style = "Heading3"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "RACSO Product Owner/Engineer"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Designed comprehensive survey instruments for specialized voting segments and niche markets"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Co-developed RACSO web application managing all aspects of survey operations from instrument design to data analysis"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Wrote RFP and analyzed bids from 1,200 vendors for research platform development"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Research Director - PCCC (Washington, DC) | 2010 - 2012"
$cur.Style = "Heading3"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Political Research & Data Analysis (FLEEM System)"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of simultaneous phone calls using emulated predictive dialer for regulated political surveys"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Built comprehensive tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Software Engineer - Salsa Labs (Washington, DC) | January 2011 - August 2011"
$cur.Style = "Heading3"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Geospatial CRM Development"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands simultaneously"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "• Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill"
$cur.Style = "Normal"
$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Text = "Programmer - Lake Research Partners (Washington, DC) | April 2008 - December 2008"
$cur.Style = "Heading3"

# Replace before[6:6] text/style in place (1 paragraph(s))
$d.Paragraphs.Item(6).Range.Text = ""

# Replace before[4:4] text/style in place (1 paragraph(s))
$d.Paragraphs.Item(4).Range.Text = "Senior data scientist and software engineer specializing in geospatial machine learning and large-scale demographic analysis. Developed algorithms that improved demographic classification accuracy from 23% to 64%, processed data across 178,000+ precincts, and built platforms serving thousands of analysts nationwide."

# Delete before[2:2] (no replacement)
$r = $d.Range($d.Paragraphs.Item(2).Range.Start, $d.Paragraphs.Item(2).Range.End)
$r.Delete()
